$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.021590333333334
$ws.Range("H2").Value = 6.064771
$ws.Range("I2").Value = 0.01116262347650641
$ws.Range("J2").Value = 0.01116262347650641
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 68.87709143691878
$ws.Range("R2").Value = 619.8938229322689
$ws.Range("S2").Value = 0.005772711098139562
$ws.Range("T2").Value = 0.005772711098139563
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.021590333333334
$ws.Range("H3").Value = 6.064771
$ws.Range("I3").Value = 0.01116262347650641
$ws.Range("J3").Value = 0.01116262347650641
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.05649099999999
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 55.968733108729
$ws.Range("R3").Value = 503.718597978561
$ws.Range("S3").Value = 0.004690838710305806
$ws.Range("T3").Value = 0.004690838710305806
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.021590333333334
$ws.Range("H4").Value = 6.064771
$ws.Range("I4").Value = 0.01116262347650641
$ws.Range("J4").Value = 0.01116262347650641
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 8.340996134675889
$ws.Range("R4").Value = 75.068965212083
$ws.Range("S4").Value = 0.0006990736680610436
$ws.Range("T4").Value = 0.0006990736680610436
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1388191398995883
$ws.Range("J5").Value = 0.1388191398995883
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 856.5601636731752
$ws.Range("R5").Value = 7709.041473058577
$ws.Range("S5").Value = 0.07178982532369227
$ws.Range("T5").Value = 0.07178982532369227
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1388191398995883
$ws.Range("J6").Value = 0.1388191398995883
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.05649099999999
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 696.0309471850981
$ws.Range("R6").Value = 6264.278524665881
$ws.Range("S6").Value = 0.05833558719801471
$ws.Range("T6").Value = 0.05833558719801471
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1388191398995883
$ws.Range("J7").Value = 0.1388191398995883
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 103.7291915971607
$ws.Range("R7").Value = 933.562724374446
$ws.Range("S7").Value = 0.008693727377881303
$ws.Range("T7").Value = 0.008693727377881301
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.823868
$ws.Range("I8").Value = 0.8500182366239053
$ws.Range("J8").Value = 0.8500182366239052
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 5244.894619102271
$ws.Range("R8").Value = 47204.05157192044
$ws.Range("S8").Value = 0.4395839130924054
$ws.Range("T8").Value = 0.4395839130924054
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.823868
$ws.Range("I9").Value = 0.8500182366239053
$ws.Range("J9").Value = 0.8500182366239052
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.05649099999999
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 4261.941104014132
$ws.Range("R9").Value = 38357.46993612719
$ws.Range("S9").Value = 0.3572008369908045
$ws.Range("T9").Value = 0.3572008369908045
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.823868
$ws.Range("I10").Value = 0.8500182366239053
$ws.Range("J10").Value = 0.8500182366239052
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 635.1552429414182
$ws.Range("R10").Value = 5716.397186472764
$ws.Range("S10").Value = 0.05323348654069531
$ws.Range("T10").Value = 0.0532334865406953
